$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (after the header row) to hold slugified
# machine-readable identifiers for each column, so that columns can be
# related to build hierarchical SKOS concepts.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "edad-grupos-quinquenales"
$ws.Range("B2").Value = "mes-codigo"
$ws.Range("C2").Value = "ccaa-nombre"
$ws.Range("D2").Value = "numero-de-contratos"
$ws.Range("E2").Value = "provincia-codigo"
$ws.Range("F2").Value = "provincia-nombre"
$ws.Range("G2").Value = "sexo"
$ws.Range("H2").Value = "mes-y-ano"
